$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add P1=14, Q1=15, matching the style (s="1") used by B1:O1 ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$excel.CutCopyMode = 0

# --- Rows 2-25: update I/K/M/O values and add new P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1

    $arr = New-Object 'object[,]' 1,2
    $arr[0,0] = 2
    $arr[0,1] = 2
    $ws.Range("P$r`:Q$r").Value = $arr
}
